$d = $word.ActiveDocument

# Locate the target paragraph robustly by matching on distinctive text
# rather than a hardcoded index, in case paragraph numbering shifts.
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Organisation/Impressum einbinden*") {
        $target = $para
    }
}
if ($target -eq $null) {
    Write-Host "ERROR: target paragraph not found"
}
$p = $target

# ---- Step 1: "Kollau: 5h – Organisation/Impressum einbinden"
#              -> "Kollau: 5h – Organisation/ "
# (the old run containing "Impressum einbinden" and the bookmark that used
#  to sit between "Organisation/" and it are removed as part of this)
$rng1 = $p.Range
$res1 = $rng1.Find.Execute("Organisation/Impressum einbinden", $true, $false, $false, $false, $false, $true, 1, $false, "Organisation/ ", 2)
Write-Host "Step1:" $res1

# ---- Step 2: remove the "Itemgrafiken" that follows "Frühwirth: 5h – "
$rng2a = $p.Range
$res2a = $rng2a.Find.Execute("Frühwirth: 5h – ")
Write-Host "Step2a:" $res2a
$rng2b = $d.Range($rng2a.End, $p.Range.End)
$res2b = $rng2b.Find.Execute("Itemgrafiken")
Write-Host "Step2b:" $res2b
$rng2b.Delete()

# ---- Step 3: remove "/PR-Folder" that follows "Zusammenstellen des Impressums"
$rng3 = $p.Range
$res3 = $rng3.Find.Execute("Zusammenstellen des Impressums/PR-Folder", $true, $false, $false, $false, $false, $true, 1, $false, "Zusammenstellen des Impressums", 2)
Write-Host "Step3:" $res3

# ---- Step 4: move the "_GoBack" bookmark so it sits right after "Easteregg"
$rng4 = $p.Range
$res4 = $rng4.Find.Execute("Easteregg")
Write-Host "Step4:" $res4
$bmRng = $rng4.Duplicate
$bmRng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRng)

Write-Host "FINAL TEXT: [" $p.Range.Text "]"
